{"js": "// Replace the \"three-digit \u00f7 one-digit\" answer cells with the new values.\n// Each lookup text is unique within the document, so a plain search +\n// whole-range replace is safe and keeps the original run formatting\n// (font/size) untouched since we only rewrite the text of the found range.\nconst replacements = [\n  [\"681\u00f73=227, 0\", \"764\u00f76=127, 2\"],\n  [\"281\u00f73=93, 2\", \"211\u00f77=30, 1\"],\n  [\"910\u00f75=182, 0\", \"693\u00f79=77, 0\"],\n  [\"766\u00f78=95, 6\", \"532\u00f78=66, 4\"],\n  [\"761\u00f77=108, 5\", \"435\u00f77=62, 1\"],\n  [\"827\u00f74=206, 3\", \"841\u00f77=120, 1\"],\n  [\"630\u00f72=315, 0\", \"489\u00f77=69, 6\"],\n  [\"353\u00f76=58, 5\", \"320\u00f78=40, 0\"],\n  [\"656\u00f79=72, 8\", \"676\u00f78=84, 4\"],\n  [\"507\u00f77=72, 3\", \"800\u00f72=400, 0\"],\n  [\"119\u00f72=59, 1\", \"536\u00f79=59, 5\"],\n  [\"230\u00f76=38, 2\", \"347\u00f78=43, 3\"],\n  [\"161\u00f76=26, 5\", \"147\u00f76=24, 3\"],\n  [\"692\u00f78=86, 4\", \"739\u00f74=184, 3\"],\n  [\"544\u00f79=60, 4\", \"850\u00f77=121, 3\"],\n  [\"594\u00f78=74, 2\", \"255\u00f76=42, 3\"],\n  [\"621\u00f77=88, 5\", \"955\u00f75=191, 0\"],\n  [\"875\u00f79=97, 2\", \"761\u00f72=380, 1\"],\n  [\"419\u00f72=209, 1\", \"800\u00f74=200, 0\"],\n  [\"872\u00f74=218, 0\", \"512\u00f76=85, 2\"],\n  [\"700\u00f73=233, 1\", \"551\u00f79=61, 2\"],\n  [\"325\u00f77=46, 3\", \"976\u00f79=108, 4\"],\n  [\"362\u00f77=51, 5\", \"837\u00f76=139, 3\"],\n  [\"374\u00f79=41, 5\", \"332\u00f77=47, 3\"],\n  [\"438\u00f73=146, 0\", \"450\u00f79=50, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the \"three-digit \u00f7 one-digit\" answer cells with the new values.\n# Each lookup text is unique within the document, so Find/Replace against\n# the whole document Content range is safe and preserves the original run\n# formatting (font/size) of the matched text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"681\u00f73=227, 0\", \"764\u00f76=127, 2\"),\n    @(\"281\u00f73=93, 2\", \"211\u00f77=30, 1\"),\n    @(\"910\u00f75=182, 0\", \"693\u00f79=77, 0\"),\n    @(\"766\u00f78=95, 6\", \"532\u00f78=66, 4\"),\n    @(\"761\u00f77=108, 5\", \"435\u00f77=62, 1\"),\n    @(\"827\u00f74=206, 3\", \"841\u00f77=120, 1\"),\n    @(\"630\u00f72=315, 0\", \"489\u00f77=69, 6\"),\n    @(\"353\u00f76=58, 5\", \"320\u00f78=40, 0\"),\n    @(\"656\u00f79=72, 8\", \"676\u00f78=84, 4\"),\n    @(\"507\u00f77=72, 3\", \"800\u00f72=400, 0\"),\n    @(\"119\u00f72=59, 1\", \"536\u00f79=59, 5\"),\n    @(\"230\u00f76=38, 2\", \"347\u00f78=43, 3\"),\n    @(\"161\u00f76=26, 5\", \"147\u00f76=24, 3\"),\n    @(\"692\u00f78=86, 4\", \"739\u00f74=184, 3\"),\n    @(\"544\u00f79=60, 4\", \"850\u00f77=121, 3\"),\n    @(\"594\u00f78=74, 2\", \"255\u00f76=42, 3\"),\n    @(\"621\u00f77=88, 5\", \"955\u00f75=191, 0\"),\n    @(\"875\u00f79=97, 2\", \"761\u00f72=380, 1\"),\n    @(\"419\u00f72=209, 1\", \"800\u00f74=200, 0\"),\n    @(\"872\u00f74=218, 0\", \"512\u00f76=85, 2\"),\n    @(\"700\u00f73=233, 1\", \"551\u00f79=61, 2\"),\n    @(\"325\u00f77=46, 3\", \"976\u00f79=108, 4\"),\n    @(\"362\u00f77=51, 5\", \"837\u00f76=139, 3\"),\n    @(\"374\u00f79=41, 5\", \"332\u00f77=47, 3\"),\n    @(\"438\u00f73=146, 0\", \"450\u00f79=50, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
